$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New statistic values after filtering save games (regenerated s_vals),
# columns B:E and the derived sum in G for rows 2-15 (F "Win" is unchanged).
$data = @{
    2  = @(0.6753301551942219, 1.667794583268128, 3.900430680208489,  0.496779210170732, 6.740334628841572)
    3  = @(3.230985683306322,  1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4  = @(1.459612070389937,  1.667794583268128, 3.900430680208489,  0.496779210170732, 7.524616544037286)
    5  = @(3.230985683306322,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    6  = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    7  = @(0.01514828764759746,1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.337247374063419)
    8  = @(3.230985683306322,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    9  = @(1.459612070389937,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    10 = @(1.459612070389937,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    11 = @(3.230985683306322,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    12 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    13 = @(3.230985683306322,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    14 = @(0.04763786555579896,0.04240448674262143,3.900430680208489, 0.496779210170732, 4.487252242677641)
    15 = @(3.230985683306322,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}
